$p = $ppt.ActivePresentation

# --- Slide 1: replace the free-form "Rectangle 1" shape with the
#     Title Slide layout's native Title/Subtitle placeholders. ---
$s1 = $p.Slides.Item(1)

# Drop the old rectangle shape.
$s1.Shapes.Item(1).Delete()

# Re-apply the "Title Slide" custom layout so the slide grows the
# ctrTitle / subTitle placeholders defined on that layout.
$titleLayout = $p.SlideMaster.CustomLayouts.Item(1)
$s1.CustomLayout = $titleLayout

# Fill in the title text; leave the subtitle placeholder empty.
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Lje;eftn4bit;vi;v"

# --- Remove the second slide (penguin picture) entirely. ---
$p.Slides.Item(2).Delete()
